$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 689, shifting existing rows (689..786) down to (690..787)
$ws.Rows.Item(689).Insert()

# Populate the newly inserted row 689 with the new data point
$ws.Cells.Item(689, 1).Value = 3
$ws.Cells.Item(689, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(689, 3).Value = "Coquimbo"
$ws.Cells.Item(689, 4).Value = 45127
$ws.Cells.Item(689, 4).NumberFormat = $ws.Cells.Item(690, 4).NumberFormat
$ws.Cells.Item(689, 5).Value = 5
$ws.Cells.Item(689, 6).Value = 100112037
$ws.Cells.Item(689, 7).Value = "Cebollín"
$ws.Cells.Item(689, 8).Value = "Sin especificar"
$ws.Cells.Item(689, 9).Value = "Primera"
$ws.Cells.Item(689, 10).Value = 240
$ws.Cells.Item(689, 11).Value = 4000
$ws.Cells.Item(689, 12).Value = 4500
$ws.Cells.Item(689, 13).Value = 4229
$ws.Cells.Item(689, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(689, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(689, 16).Value = 117
$ws.Cells.Item(689, 17).Value = 36
$ws.Cells.Item(689, 18).Value = "Hortaliza"
